# Update handback-status.xlsx timestamps for the "Generate Report for Handback" run.
# Row 5 on the "zh-cn" and "de-de" sheets corresponds to the
# 8a3b8d57-58a1-4d3f-a20c-265ee50b91c1 handoff/handback pair; its
# "Correspond Handoff Datetime" (col D) and "Correspond Handback DateTime"
# (col G) values are refreshed to the new report-generation timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-24 09:59:45"
$wsZhCn.Range("G5").Value = "2016-02-24 10:00:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-24 09:59:57"
$wsDeDe.Range("G5").Value = "2016-02-24 10:01:12"
